$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 100, shifting existing rows 100-103 down to 101-104.
$ws.Rows.Item(100).Insert()

# Populate the newly inserted row 100 with the new record.
$ws.Cells.Item(100, 1).Value = 1
$ws.Cells.Item(100, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(100, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(100, 4).Value = 45075
$ws.Cells.Item(100, 5).Value = 15
$ws.Cells.Item(100, 6).Value = 100112012
$ws.Cells.Item(100, 7).Value = "Espinaca"
$ws.Cells.Item(100, 8).Value = "Sin especificar"
$ws.Cells.Item(100, 9).Value = "Primera"
$ws.Cells.Item(100, 10).Value = 300
$ws.Cells.Item(100, 11).Value = 3000
$ws.Cells.Item(100, 12).Value = 3500
$ws.Cells.Item(100, 13).Value = 3250
$ws.Cells.Item(100, 14).Value = "$/atado 2,5 a 3 kilos"
$ws.Cells.Item(100, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(100, 16).Value = 1083
$ws.Cells.Item(100, 17).Value = 3
$ws.Cells.Item(100, 18).Value = "Hortaliza"
